# Auto-generated edit script: updates Leve profit-calculation cells (H:N)
# across all 8 sheets per the scheduled market-price data refresh.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 4023.6
$ws.Range("I8").Value = 2807.5
$ws.Range("J8").Value = 8888
$ws.Range("K8").Value = 8422.5
$ws.Range("L8").Value = 26664
$ws.Range("M8").Value = -8283.5
$ws.Range("N8").Value = -26942
$ws.Range("H17").Value = 1099.9323
$ws.Range("J17").Value = 1196.151
$ws.Range("L17").Value = 3588.453
$ws.Range("N17").Value = -3924.453
$ws.Range("H28").Value = 1010831.44
$ws.Range("I28").Value = 1389218.2
$ws.Range("J28").Value = 1800
$ws.Range("K28").Value = 1389218.2
$ws.Range("L28").Value = 1800
$ws.Range("M28").Value = -1388733.2
$ws.Range("N28").Value = -2770
$ws.Range("H62").Value = 15493.235
$ws.Range("I62").Value = 10888.5
$ws.Range("J62").Value = 22071.428
$ws.Range("K62").Value = 10888.5
$ws.Range("L62").Value = 22071.428
$ws.Range("M62").Value = -10264.5
$ws.Range("N62").Value = -23319.428
$ws.Range("H65").Value = 15493.235
$ws.Range("I65").Value = 10888.5
$ws.Range("J65").Value = 22071.428
$ws.Range("K65").Value = 54442.5
$ws.Range("L65").Value = 110357.14
$ws.Range("M65").Value = -51322.5
$ws.Range("N65").Value = -116597.14
$ws.Range("H94").Value = 19118.166
$ws.Range("I94").Value = 19118.166
$ws.Range("K94").Value = 19118.166
$ws.Range("M94").Value = -18667.166
$ws.Range("H125").Value = 22425028
$ws.Range("I125").Value = 10032
$ws.Range("K125").Value = 90288
$ws.Range("M125").Value = -87828
$ws.Range("H133").Value = 18750
$ws.Range("J133").Value = 18750
$ws.Range("L133").Value = 18750
$ws.Range("N133").Value = -28870
$ws.Range("H135").Value = 1062.6111
$ws.Range("I135").Value = 967.6042
$ws.Range("J135").Value = 1822.6666
$ws.Range("K135").Value = 8708.4378
$ws.Range("L135").Value = 16403.9994
$ws.Range("M135").Value = -6173.4378
$ws.Range("N135").Value = -21473.9994
$ws.Range("H137").Value = 47620468
$ws.Range("I137").Value = 62501164
$ws.Range("K137").Value = 187503492
$ws.Range("M137").Value = -187500942

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 25153
$ws.Range("I32").Value = 2949.3096
$ws.Range("J32").Value = 336004.66
$ws.Range("K32").Value = 2949.3096
$ws.Range("L32").Value = 336004.66
$ws.Range("M32").Value = -2662.3096
$ws.Range("N32").Value = -336578.66
$ws.Range("H61").Value = 2088.7446
$ws.Range("I61").Value = 1369.3636
$ws.Range("K61").Value = 1369.3636
$ws.Range("M61").Value = -1157.3636
$ws.Range("I102").Value = 1999.5
$ws.Range("K102").Value = 1999.5
$ws.Range("M102").Value = -377.5
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("H122").Value = 1900.8485
$ws.Range("I122").Value = 1495.238
$ws.Range("J122").Value = 2610.6667
$ws.Range("K122").Value = 4485.714
$ws.Range("L122").Value = 7832.000100000001
$ws.Range("M122").Value = -2035.714
$ws.Range("N122").Value = -12732.0001
$ws.Range("H136").Value = 2088.7446
$ws.Range("I136").Value = 1369.3636
$ws.Range("K136").Value = 4108.0908
$ws.Range("M136").Value = -1558.0908

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 39100
$ws.Range("J59").Value = 48650
$ws.Range("L59").Value = 48650
$ws.Range("N59").Value = -50344
$ws.Range("H134").Value = 25003410
$ws.Range("I134").Value = 45456804
$ws.Range("J134").Value = 4816.722
$ws.Range("K134").Value = 136370412
$ws.Range("L134").Value = 14450.166
$ws.Range("M134").Value = -136367877
$ws.Range("N134").Value = -19520.166

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1696.2941
$ws.Range("I31").Value = 1043.9131
$ws.Range("J31").Value = 3060.3635
$ws.Range("K31").Value = 1043.9131
$ws.Range("L31").Value = 3060.3635
$ws.Range("M31").Value = -748.9131
$ws.Range("N31").Value = -3650.3635
$ws.Range("H34").Value = 1696.2941
$ws.Range("I34").Value = 1043.9131
$ws.Range("J34").Value = 3060.3635
$ws.Range("K34").Value = 1043.9131
$ws.Range("L34").Value = 3060.3635
$ws.Range("M34").Value = -841.9131
$ws.Range("N34").Value = -3464.3635
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H58").Value = 1273.6
$ws.Range("I58").Value = 419.08334
$ws.Range("J58").Value = 2892.6843
$ws.Range("K58").Value = 419.08334
$ws.Range("L58").Value = 2892.6843
$ws.Range("M58").Value = -216.08334
$ws.Range("N58").Value = -3298.6843
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H122").Value = 1690.1852
$ws.Range("I122").Value = 1031.0714
$ws.Range("J122").Value = 2400
$ws.Range("K122").Value = 3093.2142
$ws.Range("L122").Value = 7200
$ws.Range("M122").Value = -643.2142000000003
$ws.Range("N122").Value = -12100
$ws.Range("H136").Value = 1273.6
$ws.Range("I136").Value = 419.08334
$ws.Range("J136").Value = 2892.6843
$ws.Range("K136").Value = 1257.25002
$ws.Range("L136").Value = 8678.052899999999
$ws.Range("M136").Value = 1292.74998
$ws.Range("N136").Value = -13778.0529

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1212.0646
$ws.Range("I5").Value = 531.64703
$ws.Range("K5").Value = 1594.94109
$ws.Range("M5").Value = -1482.94109
$ws.Range("H132").Value = 15152718
$ws.Range("I132").Value = 456
$ws.Range("J132").Value = 23811154
$ws.Range("K132").Value = 4104
$ws.Range("L132").Value = 214300386
$ws.Range("M132").Value = -1574
$ws.Range("N132").Value = -214305446
$ws.Range("H135").Value = 1212.0646
$ws.Range("I135").Value = 531.64703
$ws.Range("K135").Value = 4784.82327
$ws.Range("M135").Value = -2249.82327

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 7156.5557
$ws.Range("I17").Value = 12875
$ws.Range("J17").Value = 2581.8
$ws.Range("K17").Value = 12875
$ws.Range("L17").Value = 2581.8
$ws.Range("M17").Value = -12707
$ws.Range("N17").Value = -2917.8
$ws.Range("H70").Value = 6660.826
$ws.Range("I70").Value = 7031.1875
$ws.Range("J70").Value = 5814.2856
$ws.Range("K70").Value = 7031.1875
$ws.Range("L70").Value = 5814.2856
$ws.Range("M70").Value = -6761.1875
$ws.Range("N70").Value = -6354.2856
$ws.Range("H73").Value = 6660.826
$ws.Range("I73").Value = 7031.1875
$ws.Range("J73").Value = 5814.2856
$ws.Range("K73").Value = 7031.1875
$ws.Range("L73").Value = 5814.2856
$ws.Range("M73").Value = -6095.1875
$ws.Range("N73").Value = -7686.2856
$ws.Range("H102").Value = 2354.7666
$ws.Range("I102").Value = 1862.7826
$ws.Range("J102").Value = 3971.2856
$ws.Range("K102").Value = 1862.7826
$ws.Range("L102").Value = 3971.2856
$ws.Range("M102").Value = -240.7826
$ws.Range("N102").Value = -7215.2856

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2908.12
$ws.Range("I7").Value = 2172
$ws.Range("J7").Value = 3194.389
$ws.Range("K7").Value = 2172
$ws.Range("L7").Value = 3194.389
$ws.Range("M7").Value = -2060
$ws.Range("N7").Value = -3418.389
$ws.Range("H16").Value = 1373.75
$ws.Range("I16").Value = 1810.625
$ws.Range("K16").Value = 1810.625
$ws.Range("M16").Value = -1640.625
$ws.Range("H40").Value = 4375.2
$ws.Range("I40").Value = 3252
$ws.Range("J40").Value = 4500
$ws.Range("K40").Value = 3252
$ws.Range("L40").Value = 4500
$ws.Range("M40").Value = -3116
$ws.Range("N40").Value = -4772
$ws.Range("H55").Value = 800.3333
$ws.Range("I55").Value = 799
$ws.Range("J55").Value = 801
$ws.Range("K55").Value = 799
$ws.Range("L55").Value = 801
$ws.Range("M55").Value = -626
$ws.Range("N55").Value = -1147
$ws.Range("H61").Value = 6539.1665
$ws.Range("I61").Value = 6546.7827
$ws.Range("J61").Value = 6514.143
$ws.Range("K61").Value = 6546.7827
$ws.Range("L61").Value = 6514.143
$ws.Range("M61").Value = -6344.7827
$ws.Range("N61").Value = -6918.143
$ws.Range("H113").Value = 6539.1665
$ws.Range("I113").Value = 6546.7827
$ws.Range("J113").Value = 6514.143
$ws.Range("K113").Value = 6546.7827
$ws.Range("L113").Value = 6514.143
$ws.Range("M113").Value = -4376.7827
$ws.Range("N113").Value = -10854.143
$ws.Range("H126").Value = 2908.12
$ws.Range("I126").Value = 2172
$ws.Range("J126").Value = 3194.389
$ws.Range("K126").Value = 6516
$ws.Range("L126").Value = 9583.167000000001
$ws.Range("M126").Value = -4046
$ws.Range("N126").Value = -14523.167
$ws.Range("H136").Value = 5156.9688
$ws.Range("I136").Value = 3079.68
$ws.Range("J136").Value = 12575.857
$ws.Range("K136").Value = 9239.039999999999
$ws.Range("L136").Value = 37727.571
$ws.Range("M136").Value = -6689.039999999999
$ws.Range("N136").Value = -42827.571

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 247.63158
$ws.Range("I113").Value = 246.57143
$ws.Range("J113").Value = 250.6
$ws.Range("K113").Value = 739.71429
$ws.Range("L113").Value = 751.8
$ws.Range("M113").Value = 1430.28571
$ws.Range("N113").Value = -5091.8
$ws.Range("H132").Value = 21742200
$ws.Range("I132").Value = 31253162
$ws.Range("K132").Value = 93759486
$ws.Range("M132").Value = -93756956
$ws.Range("H136").Value = 10449569
$ws.Range("I136").Value = 13931661
$ws.Range("K136").Value = 41794983
$ws.Range("M136").Value = -41792433
